# Update "想去人数" (interested-count) values in column F on the
# "展览" (Exhibition) and "全部类型" (All types) sheets to reflect the
# latest scrape, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row -> new F value, for worksheet "展览"
$updatesExhibition = @{
    2  = 37
    5  = 187
    6  = 3804
    11 = 84
    13 = 171
    14 = 934
    21 = 3336
    22 = 5677
    26 = 512
    27 = 41
    28 = 3323
    29 = 345
    30 = 15
    31 = 2429
    33 = 513
    36 = 253
    39 = 1003
    40 = 885
    41 = 12
}

# Row -> new F value, for worksheet "全部类型"
$updatesAllTypes = @{
    2  = 37
    5  = 187
    6  = 3804
    12 = 84
    14 = 171
    15 = 934
    22 = 3336
    23 = 5677
    27 = 512
    28 = 41
    29 = 3323
    30 = 345
    31 = 15
    32 = 2429
    34 = 513
    37 = 253
    40 = 1003
    41 = 885
    42 = 12
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $updatesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAllTypes.Keys) {
    $wsAllTypes.Range("F$row").Value = $updatesAllTypes[$row]
}
